$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 675
$ws.Range("B2").Value = 45767.22928240741
$ws.Range("E2").Value = 60
$ws.Range("F2").Value = 0.9
$ws.Range("G2").Value = 0.9
$ws.Range("H2").Value = 1.8
$ws.Range("I2").Value = "Chicken Wrap (x1)"
